$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension to reflect two newly added rows (A1:R438)

$ws.Cells.Item(375,4).Value = 44505   # D375 Fecha
$ws.Cells.Item(375,10).Value = 3400  # J375 Volumen
$ws.Cells.Item(375,11).Value = 600  # K375 Precio minimo
$ws.Cells.Item(375,12).Value = 700  # L375 Precio maximo
$ws.Cells.Item(375,13).Value = 650  # M375 Precio promedio ponderado
$ws.Cells.Item(375,16).Value = 650  # P375 Precio $/Kg

$ws.Cells.Item(376,4).Value = 44505   # D376 Fecha
$ws.Cells.Item(376,10).Value = 1800  # J376 Volumen
$ws.Cells.Item(376,11).Value = 500  # K376 Precio minimo
$ws.Cells.Item(376,12).Value = 550  # L376 Precio maximo
$ws.Cells.Item(376,13).Value = 525  # M376 Precio promedio ponderado
$ws.Cells.Item(376,16).Value = 525  # P376 Precio $/Kg

$ws.Cells.Item(377,4).Value = 44487   # D377 Fecha
$ws.Cells.Item(377,10).Value = 3000  # J377 Volumen
$ws.Cells.Item(377,11).Value = 600  # K377 Precio minimo
$ws.Cells.Item(377,12).Value = 700  # L377 Precio maximo
$ws.Cells.Item(377,13).Value = 650  # M377 Precio promedio ponderado
$ws.Cells.Item(377,16).Value = 650  # P377 Precio $/Kg

$ws.Cells.Item(378,4).Value = 44487   # D378 Fecha
$ws.Cells.Item(378,10).Value = 1800  # J378 Volumen
$ws.Cells.Item(378,11).Value = 500  # K378 Precio minimo
$ws.Cells.Item(378,12).Value = 550  # L378 Precio maximo
$ws.Cells.Item(378,13).Value = 525  # M378 Precio promedio ponderado
$ws.Cells.Item(378,16).Value = 525  # P378 Precio $/Kg

$ws.Cells.Item(379,4).Value = 44425   # D379 Fecha
$ws.Cells.Item(379,10).Value = 2200  # J379 Volumen
$ws.Cells.Item(379,11).Value = 650  # K379 Precio minimo
$ws.Cells.Item(379,12).Value = 700  # L379 Precio maximo
$ws.Cells.Item(379,13).Value = 675  # M379 Precio promedio ponderado
$ws.Cells.Item(379,16).Value = 675  # P379 Precio $/Kg

$ws.Cells.Item(380,4).Value = 44425   # D380 Fecha
$ws.Cells.Item(380,10).Value = 1400  # J380 Volumen
$ws.Cells.Item(380,11).Value = 550  # K380 Precio minimo
$ws.Cells.Item(380,12).Value = 600  # L380 Precio maximo
$ws.Cells.Item(380,13).Value = 575  # M380 Precio promedio ponderado
$ws.Cells.Item(380,16).Value = 575  # P380 Precio $/Kg

$ws.Cells.Item(381,4).Value = 44343   # D381 Fecha
$ws.Cells.Item(381,10).Value = 2440  # J381 Volumen
$ws.Cells.Item(381,11).Value = 650  # K381 Precio minimo
$ws.Cells.Item(381,12).Value = 700  # L381 Precio maximo
$ws.Cells.Item(381,13).Value = 675  # M381 Precio promedio ponderado
$ws.Cells.Item(381,16).Value = 675  # P381 Precio $/Kg

$ws.Cells.Item(382,4).Value = 44343   # D382 Fecha
$ws.Cells.Item(382,10).Value = 1560  # J382 Volumen
$ws.Cells.Item(382,11).Value = 500  # K382 Precio minimo
$ws.Cells.Item(382,12).Value = 550  # L382 Precio maximo
$ws.Cells.Item(382,13).Value = 525  # M382 Precio promedio ponderado
$ws.Cells.Item(382,16).Value = 525  # P382 Precio $/Kg

$ws.Cells.Item(383,4).Value = 44370   # D383 Fecha
$ws.Cells.Item(383,10).Value = 3200  # J383 Volumen
$ws.Cells.Item(383,11).Value = 600  # K383 Precio minimo
$ws.Cells.Item(383,12).Value = 700  # L383 Precio maximo
$ws.Cells.Item(383,13).Value = 650  # M383 Precio promedio ponderado
$ws.Cells.Item(383,16).Value = 650  # P383 Precio $/Kg

$ws.Cells.Item(384,4).Value = 44370   # D384 Fecha
$ws.Cells.Item(384,10).Value = 1680  # J384 Volumen
$ws.Cells.Item(384,11).Value = 500  # K384 Precio minimo
$ws.Cells.Item(384,12).Value = 550  # L384 Precio maximo
$ws.Cells.Item(384,13).Value = 525  # M384 Precio promedio ponderado
$ws.Cells.Item(384,16).Value = 525  # P384 Precio $/Kg

$ws.Cells.Item(385,4).Value = 44449   # D385 Fecha
$ws.Cells.Item(385,10).Value = 3600  # J385 Volumen
$ws.Cells.Item(385,11).Value = 650  # K385 Precio minimo
$ws.Cells.Item(385,12).Value = 700  # L385 Precio maximo
$ws.Cells.Item(385,13).Value = 675  # M385 Precio promedio ponderado
$ws.Cells.Item(385,16).Value = 675  # P385 Precio $/Kg

$ws.Cells.Item(386,4).Value = 44449   # D386 Fecha
$ws.Cells.Item(386,10).Value = 1900  # J386 Volumen
$ws.Cells.Item(386,11).Value = 550  # K386 Precio minimo
$ws.Cells.Item(386,12).Value = 600  # L386 Precio maximo
$ws.Cells.Item(386,13).Value = 575  # M386 Precio promedio ponderado
$ws.Cells.Item(386,16).Value = 575  # P386 Precio $/Kg

$ws.Cells.Item(387,4).Value = 44168   # D387 Fecha
$ws.Cells.Item(387,10).Value = 2000  # J387 Volumen
$ws.Cells.Item(387,11).Value = 450  # K387 Precio minimo
$ws.Cells.Item(387,12).Value = 500  # L387 Precio maximo
$ws.Cells.Item(387,13).Value = 475  # M387 Precio promedio ponderado
$ws.Cells.Item(387,16).Value = 475  # P387 Precio $/Kg

$ws.Cells.Item(388,4).Value = 44168   # D388 Fecha
$ws.Cells.Item(388,10).Value = 1300  # J388 Volumen
$ws.Cells.Item(388,11).Value = 350  # K388 Precio minimo
$ws.Cells.Item(388,12).Value = 400  # L388 Precio maximo
$ws.Cells.Item(388,13).Value = 375  # M388 Precio promedio ponderado
$ws.Cells.Item(388,16).Value = 375  # P388 Precio $/Kg

$ws.Cells.Item(389,4).Value = 44175   # D389 Fecha
$ws.Cells.Item(389,10).Value = 2000  # J389 Volumen
$ws.Cells.Item(389,11).Value = 550  # K389 Precio minimo
$ws.Cells.Item(389,12).Value = 600  # L389 Precio maximo
$ws.Cells.Item(389,13).Value = 575  # M389 Precio promedio ponderado
$ws.Cells.Item(389,16).Value = 575  # P389 Precio $/Kg

$ws.Cells.Item(390,4).Value = 44175   # D390 Fecha
$ws.Cells.Item(390,10).Value = 1400  # J390 Volumen
$ws.Cells.Item(390,11).Value = 450  # K390 Precio minimo
$ws.Cells.Item(390,12).Value = 500  # L390 Precio maximo
$ws.Cells.Item(390,13).Value = 475  # M390 Precio promedio ponderado
$ws.Cells.Item(390,16).Value = 475  # P390 Precio $/Kg

$ws.Cells.Item(391,4).Value = 44392   # D391 Fecha
$ws.Cells.Item(391,10).Value = 2300  # J391 Volumen
$ws.Cells.Item(391,11).Value = 600  # K391 Precio minimo
$ws.Cells.Item(391,12).Value = 700  # L391 Precio maximo
$ws.Cells.Item(391,13).Value = 650  # M391 Precio promedio ponderado
$ws.Cells.Item(391,16).Value = 650  # P391 Precio $/Kg

$ws.Cells.Item(392,4).Value = 44392   # D392 Fecha
$ws.Cells.Item(392,10).Value = 1400  # J392 Volumen
$ws.Cells.Item(392,11).Value = 500  # K392 Precio minimo
$ws.Cells.Item(392,12).Value = 550  # L392 Precio maximo
$ws.Cells.Item(392,13).Value = 525  # M392 Precio promedio ponderado
$ws.Cells.Item(392,16).Value = 525  # P392 Precio $/Kg

$ws.Cells.Item(393,4).Value = 44286   # D393 Fecha
$ws.Cells.Item(393,10).Value = 3200  # J393 Volumen
$ws.Cells.Item(393,11).Value = 750  # K393 Precio minimo
$ws.Cells.Item(393,12).Value = 800  # L393 Precio maximo
$ws.Cells.Item(393,13).Value = 775  # M393 Precio promedio ponderado
$ws.Cells.Item(393,16).Value = 775  # P393 Precio $/Kg

$ws.Cells.Item(394,4).Value = 44286   # D394 Fecha
$ws.Cells.Item(394,10).Value = 1720  # J394 Volumen
$ws.Cells.Item(394,11).Value = 650  # K394 Precio minimo
$ws.Cells.Item(394,12).Value = 700  # L394 Precio maximo
$ws.Cells.Item(394,13).Value = 675  # M394 Precio promedio ponderado
$ws.Cells.Item(394,16).Value = 675  # P394 Precio $/Kg

$ws.Cells.Item(395,4).Value = 44473   # D395 Fecha
$ws.Cells.Item(395,10).Value = 2600  # J395 Volumen
$ws.Cells.Item(395,11).Value = 650  # K395 Precio minimo
$ws.Cells.Item(395,12).Value = 700  # L395 Precio maximo
$ws.Cells.Item(395,13).Value = 675  # M395 Precio promedio ponderado
$ws.Cells.Item(395,16).Value = 675  # P395 Precio $/Kg

$ws.Cells.Item(396,4).Value = 44473   # D396 Fecha
$ws.Cells.Item(396,10).Value = 1600  # J396 Volumen
$ws.Cells.Item(396,11).Value = 550  # K396 Precio minimo
$ws.Cells.Item(396,12).Value = 600  # L396 Precio maximo
$ws.Cells.Item(396,13).Value = 575  # M396 Precio promedio ponderado
$ws.Cells.Item(396,16).Value = 575  # P396 Precio $/Kg

$ws.Cells.Item(397,4).Value = 44400   # D397 Fecha
$ws.Cells.Item(397,10).Value = 3600  # J397 Volumen
$ws.Cells.Item(397,11).Value = 700  # K397 Precio minimo
$ws.Cells.Item(397,12).Value = 800  # L397 Precio maximo
$ws.Cells.Item(397,13).Value = 750  # M397 Precio promedio ponderado
$ws.Cells.Item(397,16).Value = 750  # P397 Precio $/Kg

$ws.Cells.Item(398,4).Value = 44400   # D398 Fecha
$ws.Cells.Item(398,10).Value = 1800  # J398 Volumen
$ws.Cells.Item(398,11).Value = 500  # K398 Precio minimo
$ws.Cells.Item(398,12).Value = 600  # L398 Precio maximo
$ws.Cells.Item(398,13).Value = 550  # M398 Precio promedio ponderado
$ws.Cells.Item(398,16).Value = 550  # P398 Precio $/Kg

$ws.Cells.Item(399,4).Value = 44484   # D399 Fecha
$ws.Cells.Item(399,10).Value = 3600  # J399 Volumen
$ws.Cells.Item(399,11).Value = 600  # K399 Precio minimo
$ws.Cells.Item(399,12).Value = 700  # L399 Precio maximo
$ws.Cells.Item(399,13).Value = 650  # M399 Precio promedio ponderado
$ws.Cells.Item(399,16).Value = 650  # P399 Precio $/Kg

$ws.Cells.Item(400,4).Value = 44484   # D400 Fecha
$ws.Cells.Item(400,10).Value = 1800  # J400 Volumen
$ws.Cells.Item(400,11).Value = 500  # K400 Precio minimo
$ws.Cells.Item(400,12).Value = 550  # L400 Precio maximo
$ws.Cells.Item(400,13).Value = 525  # M400 Precio promedio ponderado
$ws.Cells.Item(400,16).Value = 525  # P400 Precio $/Kg

$ws.Cells.Item(401,4).Value = 44181   # D401 Fecha
$ws.Cells.Item(401,10).Value = 2000  # J401 Volumen
$ws.Cells.Item(401,11).Value = 550  # K401 Precio minimo
$ws.Cells.Item(401,12).Value = 600  # L401 Precio maximo
$ws.Cells.Item(401,13).Value = 575  # M401 Precio promedio ponderado
$ws.Cells.Item(401,16).Value = 575  # P401 Precio $/Kg

$ws.Cells.Item(402,4).Value = 44181   # D402 Fecha
$ws.Cells.Item(402,10).Value = 1200  # J402 Volumen
$ws.Cells.Item(402,11).Value = 450  # K402 Precio minimo
$ws.Cells.Item(402,12).Value = 500  # L402 Precio maximo
$ws.Cells.Item(402,13).Value = 475  # M402 Precio promedio ponderado
$ws.Cells.Item(402,16).Value = 475  # P402 Precio $/Kg

$ws.Cells.Item(403,4).Value = 44494   # D403 Fecha
$ws.Cells.Item(403,10).Value = 2700  # J403 Volumen
$ws.Cells.Item(403,11).Value = 600  # K403 Precio minimo
$ws.Cells.Item(403,12).Value = 700  # L403 Precio maximo
$ws.Cells.Item(403,13).Value = 650  # M403 Precio promedio ponderado
$ws.Cells.Item(403,16).Value = 650  # P403 Precio $/Kg

$ws.Cells.Item(404,4).Value = 44494   # D404 Fecha
$ws.Cells.Item(404,10).Value = 1600  # J404 Volumen
$ws.Cells.Item(404,11).Value = 500  # K404 Precio minimo
$ws.Cells.Item(404,12).Value = 550  # L404 Precio maximo
$ws.Cells.Item(404,13).Value = 525  # M404 Precio promedio ponderado
$ws.Cells.Item(404,16).Value = 525  # P404 Precio $/Kg

$ws.Cells.Item(405,4).Value = 44342   # D405 Fecha
$ws.Cells.Item(405,10).Value = 3260  # J405 Volumen
$ws.Cells.Item(405,11).Value = 650  # K405 Precio minimo
$ws.Cells.Item(405,12).Value = 700  # L405 Precio maximo
$ws.Cells.Item(405,13).Value = 675  # M405 Precio promedio ponderado
$ws.Cells.Item(405,16).Value = 675  # P405 Precio $/Kg

$ws.Cells.Item(406,4).Value = 44342   # D406 Fecha
$ws.Cells.Item(406,10).Value = 1680  # J406 Volumen
$ws.Cells.Item(406,11).Value = 500  # K406 Precio minimo
$ws.Cells.Item(406,12).Value = 550  # L406 Precio maximo
$ws.Cells.Item(406,13).Value = 525  # M406 Precio promedio ponderado
$ws.Cells.Item(406,16).Value = 525  # P406 Precio $/Kg

$ws.Cells.Item(407,4).Value = 44445   # D407 Fecha
$ws.Cells.Item(407,10).Value = 3340  # J407 Volumen
$ws.Cells.Item(407,11).Value = 650  # K407 Precio minimo
$ws.Cells.Item(407,12).Value = 700  # L407 Precio maximo
$ws.Cells.Item(407,13).Value = 675  # M407 Precio promedio ponderado
$ws.Cells.Item(407,16).Value = 675  # P407 Precio $/Kg

$ws.Cells.Item(408,4).Value = 44445   # D408 Fecha
$ws.Cells.Item(408,10).Value = 1680  # J408 Volumen
$ws.Cells.Item(408,11).Value = 550  # K408 Precio minimo
$ws.Cells.Item(408,12).Value = 600  # L408 Precio maximo
$ws.Cells.Item(408,13).Value = 575  # M408 Precio promedio ponderado
$ws.Cells.Item(408,16).Value = 575  # P408 Precio $/Kg

$ws.Cells.Item(409,4).Value = 44328   # D409 Fecha
$ws.Cells.Item(409,10).Value = 3260  # J409 Volumen
$ws.Cells.Item(409,11).Value = 650  # K409 Precio minimo
$ws.Cells.Item(409,12).Value = 700  # L409 Precio maximo
$ws.Cells.Item(409,13).Value = 675  # M409 Precio promedio ponderado
$ws.Cells.Item(409,16).Value = 675  # P409 Precio $/Kg

$ws.Cells.Item(410,4).Value = 44328   # D410 Fecha
$ws.Cells.Item(410,10).Value = 1600  # J410 Volumen
$ws.Cells.Item(410,11).Value = 500  # K410 Precio minimo
$ws.Cells.Item(410,12).Value = 550  # L410 Precio maximo
$ws.Cells.Item(410,13).Value = 525  # M410 Precio promedio ponderado
$ws.Cells.Item(410,16).Value = 525  # P410 Precio $/Kg

$ws.Cells.Item(411,4).Value = 44301   # D411 Fecha
$ws.Cells.Item(411,10).Value = 2500  # J411 Volumen
$ws.Cells.Item(411,11).Value = 650  # K411 Precio minimo
$ws.Cells.Item(411,12).Value = 700  # L411 Precio maximo
$ws.Cells.Item(411,13).Value = 675  # M411 Precio promedio ponderado
$ws.Cells.Item(411,16).Value = 675  # P411 Precio $/Kg

$ws.Cells.Item(412,4).Value = 44301   # D412 Fecha
$ws.Cells.Item(412,10).Value = 1400  # J412 Volumen
$ws.Cells.Item(412,11).Value = 550  # K412 Precio minimo
$ws.Cells.Item(412,12).Value = 600  # L412 Precio maximo
$ws.Cells.Item(412,13).Value = 575  # M412 Precio promedio ponderado
$ws.Cells.Item(412,16).Value = 575  # P412 Precio $/Kg

$ws.Cells.Item(413,4).Value = 44330   # D413 Fecha
$ws.Cells.Item(413,10).Value = 3300  # J413 Volumen
$ws.Cells.Item(413,11).Value = 650  # K413 Precio minimo
$ws.Cells.Item(413,12).Value = 700  # L413 Precio maximo
$ws.Cells.Item(413,13).Value = 675  # M413 Precio promedio ponderado
$ws.Cells.Item(413,16).Value = 675  # P413 Precio $/Kg

$ws.Cells.Item(414,4).Value = 44330   # D414 Fecha
$ws.Cells.Item(414,10).Value = 1660  # J414 Volumen
$ws.Cells.Item(414,11).Value = 500  # K414 Precio minimo
$ws.Cells.Item(414,12).Value = 550  # L414 Precio maximo
$ws.Cells.Item(414,13).Value = 525  # M414 Precio promedio ponderado
$ws.Cells.Item(414,16).Value = 525  # P414 Precio $/Kg

$ws.Cells.Item(415,4).Value = 44270   # D415 Fecha
$ws.Cells.Item(415,10).Value = 2800  # J415 Volumen
$ws.Cells.Item(415,11).Value = 850  # K415 Precio minimo
$ws.Cells.Item(415,12).Value = 900  # L415 Precio maximo
$ws.Cells.Item(415,13).Value = 875  # M415 Precio promedio ponderado
$ws.Cells.Item(415,16).Value = 875  # P415 Precio $/Kg

$ws.Cells.Item(416,4).Value = 44270   # D416 Fecha
$ws.Cells.Item(416,10).Value = 1540  # J416 Volumen
$ws.Cells.Item(416,11).Value = 750  # K416 Precio minimo
$ws.Cells.Item(416,12).Value = 800  # L416 Precio maximo
$ws.Cells.Item(416,13).Value = 775  # M416 Precio promedio ponderado
$ws.Cells.Item(416,16).Value = 775  # P416 Precio $/Kg

$ws.Cells.Item(417,4).Value = 44295   # D417 Fecha
$ws.Cells.Item(417,10).Value = 3200  # J417 Volumen
$ws.Cells.Item(417,11).Value = 650  # K417 Precio minimo
$ws.Cells.Item(417,12).Value = 700  # L417 Precio maximo
$ws.Cells.Item(417,13).Value = 675  # M417 Precio promedio ponderado
$ws.Cells.Item(417,16).Value = 675  # P417 Precio $/Kg

$ws.Cells.Item(418,4).Value = 44295   # D418 Fecha
$ws.Cells.Item(418,10).Value = 1660  # J418 Volumen
$ws.Cells.Item(418,11).Value = 550  # K418 Precio minimo
$ws.Cells.Item(418,12).Value = 600  # L418 Precio maximo
$ws.Cells.Item(418,13).Value = 575  # M418 Precio promedio ponderado
$ws.Cells.Item(418,16).Value = 575  # P418 Precio $/Kg

$ws.Cells.Item(419,4).Value = 44217   # D419 Fecha
$ws.Cells.Item(419,10).Value = 2000  # J419 Volumen
$ws.Cells.Item(419,11).Value = 650  # K419 Precio minimo
$ws.Cells.Item(419,12).Value = 700  # L419 Precio maximo
$ws.Cells.Item(419,13).Value = 675  # M419 Precio promedio ponderado
$ws.Cells.Item(419,16).Value = 675  # P419 Precio $/Kg

$ws.Cells.Item(420,4).Value = 44217   # D420 Fecha
$ws.Cells.Item(420,10).Value = 1440  # J420 Volumen
$ws.Cells.Item(420,11).Value = 500  # K420 Precio minimo
$ws.Cells.Item(420,12).Value = 550  # L420 Precio maximo
$ws.Cells.Item(420,13).Value = 525  # M420 Precio promedio ponderado
$ws.Cells.Item(420,16).Value = 525  # P420 Precio $/Kg

$ws.Cells.Item(421,4).Value = 44421   # D421 Fecha
$ws.Cells.Item(421,10).Value = 3600  # J421 Volumen
$ws.Cells.Item(421,11).Value = 650  # K421 Precio minimo
$ws.Cells.Item(421,12).Value = 700  # L421 Precio maximo
$ws.Cells.Item(421,13).Value = 675  # M421 Precio promedio ponderado
$ws.Cells.Item(421,16).Value = 675  # P421 Precio $/Kg

$ws.Cells.Item(422,4).Value = 44421   # D422 Fecha
$ws.Cells.Item(422,10).Value = 1800  # J422 Volumen
$ws.Cells.Item(422,11).Value = 550  # K422 Precio minimo
$ws.Cells.Item(422,12).Value = 600  # L422 Precio maximo
$ws.Cells.Item(422,13).Value = 575  # M422 Precio promedio ponderado
$ws.Cells.Item(422,16).Value = 575  # P422 Precio $/Kg

$ws.Cells.Item(423,4).Value = 44383   # D423 Fecha
$ws.Cells.Item(423,10).Value = 2360  # J423 Volumen
$ws.Cells.Item(423,11).Value = 600  # K423 Precio minimo
$ws.Cells.Item(423,12).Value = 700  # L423 Precio maximo
$ws.Cells.Item(423,13).Value = 650  # M423 Precio promedio ponderado
$ws.Cells.Item(423,16).Value = 650  # P423 Precio $/Kg

$ws.Cells.Item(424,4).Value = 44383   # D424 Fecha
$ws.Cells.Item(424,10).Value = 1360  # J424 Volumen
$ws.Cells.Item(424,11).Value = 500  # K424 Precio minimo
$ws.Cells.Item(424,12).Value = 550  # L424 Precio maximo
$ws.Cells.Item(424,13).Value = 525  # M424 Precio promedio ponderado
$ws.Cells.Item(424,16).Value = 525  # P424 Precio $/Kg

$ws.Cells.Item(425,4).Value = 44244   # D425 Fecha
$ws.Cells.Item(425,10).Value = 3000  # J425 Volumen
$ws.Cells.Item(425,11).Value = 750  # K425 Precio minimo
$ws.Cells.Item(425,12).Value = 800  # L425 Precio maximo
$ws.Cells.Item(425,13).Value = 775  # M425 Precio promedio ponderado
$ws.Cells.Item(425,16).Value = 775  # P425 Precio $/Kg

$ws.Cells.Item(426,4).Value = 44244   # D426 Fecha
$ws.Cells.Item(426,10).Value = 1600  # J426 Volumen
$ws.Cells.Item(426,11).Value = 650  # K426 Precio minimo
$ws.Cells.Item(426,12).Value = 700  # L426 Precio maximo
$ws.Cells.Item(426,13).Value = 675  # M426 Precio promedio ponderado
$ws.Cells.Item(426,16).Value = 675  # P426 Precio $/Kg

$ws.Cells.Item(427,4).Value = 44307   # D427 Fecha
$ws.Cells.Item(427,10).Value = 3300  # J427 Volumen
$ws.Cells.Item(427,11).Value = 650  # K427 Precio minimo
$ws.Cells.Item(427,12).Value = 700  # L427 Precio maximo
$ws.Cells.Item(427,13).Value = 675  # M427 Precio promedio ponderado
$ws.Cells.Item(427,16).Value = 675  # P427 Precio $/Kg

$ws.Cells.Item(428,4).Value = 44307   # D428 Fecha
$ws.Cells.Item(428,10).Value = 1720  # J428 Volumen
$ws.Cells.Item(428,11).Value = 550  # K428 Precio minimo
$ws.Cells.Item(428,12).Value = 600  # L428 Precio maximo
$ws.Cells.Item(428,13).Value = 575  # M428 Precio promedio ponderado
$ws.Cells.Item(428,16).Value = 575  # P428 Precio $/Kg

$ws.Cells.Item(429,4).Value = 44273   # D429 Fecha
$ws.Cells.Item(429,10).Value = 2500  # J429 Volumen
$ws.Cells.Item(429,11).Value = 750  # K429 Precio minimo
$ws.Cells.Item(429,12).Value = 800  # L429 Precio maximo
$ws.Cells.Item(429,13).Value = 775  # M429 Precio promedio ponderado
$ws.Cells.Item(429,16).Value = 775  # P429 Precio $/Kg

$ws.Cells.Item(430,4).Value = 44273   # D430 Fecha
$ws.Cells.Item(430,10).Value = 1440  # J430 Volumen
$ws.Cells.Item(430,11).Value = 650  # K430 Precio minimo
$ws.Cells.Item(430,12).Value = 700  # L430 Precio maximo
$ws.Cells.Item(430,13).Value = 675  # M430 Precio promedio ponderado
$ws.Cells.Item(430,16).Value = 675  # P430 Precio $/Kg

$ws.Cells.Item(431,4).Value = 44433   # D431 Fecha
$ws.Cells.Item(431,10).Value = 3320  # J431 Volumen
$ws.Cells.Item(431,11).Value = 650  # K431 Precio minimo
$ws.Cells.Item(431,12).Value = 700  # L431 Precio maximo
$ws.Cells.Item(431,13).Value = 675  # M431 Precio promedio ponderado
$ws.Cells.Item(431,16).Value = 675  # P431 Precio $/Kg

$ws.Cells.Item(432,4).Value = 44433   # D432 Fecha
$ws.Cells.Item(432,10).Value = 1640  # J432 Volumen
$ws.Cells.Item(432,11).Value = 550  # K432 Precio minimo
$ws.Cells.Item(432,12).Value = 600  # L432 Precio maximo
$ws.Cells.Item(432,13).Value = 575  # M432 Precio promedio ponderado
$ws.Cells.Item(432,16).Value = 575  # P432 Precio $/Kg

$ws.Cells.Item(433,4).Value = 44302   # D433 Fecha
$ws.Cells.Item(433,10).Value = 3200  # J433 Volumen
$ws.Cells.Item(433,11).Value = 650  # K433 Precio minimo
$ws.Cells.Item(433,12).Value = 700  # L433 Precio maximo
$ws.Cells.Item(433,13).Value = 675  # M433 Precio promedio ponderado
$ws.Cells.Item(433,16).Value = 675  # P433 Precio $/Kg

$ws.Cells.Item(434,4).Value = 44302   # D434 Fecha
$ws.Cells.Item(434,10).Value = 1660  # J434 Volumen
$ws.Cells.Item(434,11).Value = 550  # K434 Precio minimo
$ws.Cells.Item(434,12).Value = 600  # L434 Precio maximo
$ws.Cells.Item(434,13).Value = 575  # M434 Precio promedio ponderado
$ws.Cells.Item(434,16).Value = 575  # P434 Precio $/Kg

$ws.Cells.Item(435,4).Value = 44179   # D435 Fecha
$ws.Cells.Item(435,10).Value = 2600  # J435 Volumen
$ws.Cells.Item(435,11).Value = 550  # K435 Precio minimo
$ws.Cells.Item(435,12).Value = 600  # L435 Precio maximo
$ws.Cells.Item(435,13).Value = 575  # M435 Precio promedio ponderado
$ws.Cells.Item(435,16).Value = 575  # P435 Precio $/Kg

$ws.Cells.Item(436,4).Value = 44179   # D436 Fecha
$ws.Cells.Item(436,10).Value = 1500  # J436 Volumen
$ws.Cells.Item(436,11).Value = 450  # K436 Precio minimo
$ws.Cells.Item(436,12).Value = 500  # L436 Precio maximo
$ws.Cells.Item(436,13).Value = 475  # M436 Precio promedio ponderado
$ws.Cells.Item(436,16).Value = 475  # P436 Precio $/Kg

# New row 437
$ws.Cells.Item(437,1).Value = 8
$ws.Cells.Item(437,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(437,3).Value = "Coquimbo"
$ws.Cells.Item(437,4).Value = 44491
$ws.Cells.Item(437,5).Value = 4
$ws.Cells.Item(437,6).Value = 100112008
$ws.Cells.Item(437,7).Value = "Coliflor"
$ws.Cells.Item(437,8).Value = "Sin especificar"
$ws.Cells.Item(437,9).Value = "Primera"
$ws.Cells.Item(437,10).Value = 3600
$ws.Cells.Item(437,11).Value = 600
$ws.Cells.Item(437,12).Value = 700
$ws.Cells.Item(437,13).Value = 650
$ws.Cells.Item(437,14).Value = "$/unidad"
$ws.Cells.Item(437,15).Value = "Provincia del Elquí"
$ws.Cells.Item(437,16).Value = 650
$ws.Cells.Item(437,17).Value = 1
$ws.Cells.Item(437,18).Value = "Hortaliza"

# New row 438
$ws.Cells.Item(438,1).Value = 8
$ws.Cells.Item(438,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(438,3).Value = "Coquimbo"
$ws.Cells.Item(438,4).Value = 44491
$ws.Cells.Item(438,5).Value = 4
$ws.Cells.Item(438,6).Value = 100112008
$ws.Cells.Item(438,7).Value = "Coliflor"
$ws.Cells.Item(438,8).Value = "Sin especificar"
$ws.Cells.Item(438,9).Value = "Segunda"
$ws.Cells.Item(438,10).Value = 1800
$ws.Cells.Item(438,11).Value = 500
$ws.Cells.Item(438,12).Value = 550
$ws.Cells.Item(438,13).Value = 525
$ws.Cells.Item(438,14).Value = "$/unidad"
$ws.Cells.Item(438,15).Value = "Provincia del Elquí"
$ws.Cells.Item(438,16).Value = 525
$ws.Cells.Item(438,17).Value = 1
$ws.Cells.Item(438,18).Value = "Hortaliza"

# Apply the date number format used elsewhere in column D to the two new date cells
$ws.Cells.Item(437,4).NumberFormat = $ws.Cells.Item(435,4).NumberFormat
$ws.Cells.Item(438,4).NumberFormat = $ws.Cells.Item(435,4).NumberFormat
